$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_device")

# Append 10 new Mac-Address rows (device ids 3000166 .. 3000175) below the
# existing data, mirroring the existing rows' layout/values.
$startRow = 147
$startDeviceId = 3000166

for ($i = 0; $i -lt 10; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $startDeviceId + $i
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

$ws.Range("C152").Select() | Out-Null
